$wb = $excel.ActiveWorkbook

# New timestamp / version used across every ObjTables header string in this workbook.
$newDate = "2020-05-29 00:23:52"
$newVersion = "1.0.0"

# Map of worksheet (tab) name -> class name used inside the per-sheet header string.
$classesBySheet = [ordered]@{
    "!!Compartment"            = "Compartment"
    "!!Compound"               = "Compound"
    "!!Definition"             = "Definition"
    "!!Enzyme"                 = "Enzyme"
    "!!FbcObjective"           = "FbcObjective"
    "!!Gene"                   = "Gene"
    "!!Layout"                 = "Layout"
    "!!Measurement"            = "Measurement"
    "!!PbConfig"               = "PbConfig"
    "!!Position"               = "Position"
    "!!Protein"                = "Protein"
    "!!Quantity"               = "Quantity"
    "!!QuantityInfo"           = "QuantityInfo"
    "!!QuantityMatrix"         = "QuantityMatrix"
    "!!Reaction"                = "Reaction"
    "!!ReactionStoichiometry"  = "ReactionStoichiometry"
    "!!Regulator"              = "Regulator"
    "!!Relation"               = "Relation"
    "!!Relationship"           = "Relationship"
    "!!SparseMatrix"           = "SparseMatrix"
    "!!SparseMatrixColumn"     = "SparseMatrixColumn"
    "!!SparseMatrixOrdered"    = "SparseMatrixOrdered"
    "!!SparseMatrixRow"        = "SparseMatrixRow"
    "!!StoichiometricMatrix"   = "StoichiometricMatrix"
    "!!rxnconContingencyList"  = "rxnconContingencyList"
    "!!rxnconReactionList"     = "rxnconReactionList"
}

foreach ($sheetName in $classesBySheet.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $class = $classesBySheet[$sheetName]

    $classHeader = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' class='$class' name='$class' date='$newDate' objTablesVersion='$newVersion'"

    # Sheets are protected (no password); unprotect, edit, then restore protection.
    $ws.Unprotect()

    if ($sheetName -eq "!!Compartment") {
        # This first sheet also carries the workbook-level ObjTables banner in A1,
        # with the per-class header pushed down to A2.
        $ws.Range("A1").Value = "!!!ObjTables schema='SBtab' objTablesVersion='$newVersion' date='$newDate'"
        $ws.Range("A2").Value = $classHeader
    } else {
        $ws.Range("A1").Value = $classHeader
    }

    $ws.Protect()
}
